$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "time:timestamp"
$ws.Range("B2").Value = "datetime"

$ws.Range("A3").Value = "org:resource"
$ws.Range("B3").Value = "str"

$ws.Range("A5").Value = "SubProcessID"

$ws.Range("A6").Value = "stream:datastream"
$ws.Range("B6").Value = "dict"

$ws.Range("A7").Value = "concept:name"
